$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New TPM values per updated script output.
# Row 2 (Sending=FAPs, Target=ECs)
$ws.Range("I2").Value = 0.1309698538606213
$ws.Range("J2").Value = 0.1843806287874228
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.2317656666666667
$ws.Range("N2").Value = 0.6952970000000001
$ws.Range("O2").Value = 0.04497992035421015
$ws.Range("P2").Value = 0.04579239628237638
$ws.Range("Q2").Value = 0.07432137790311111
$ws.Range("R2").Value = 0.6688924011280001
$ws.Range("S2").Value = 0.00589101359545329
$ws.Range("T2").Value = 0.008443230820227401

# Row 3 (Sending=FAPs, Target=FAPs)
$ws.Range("I3").Value = 0.1309698538606213
$ws.Range("J3").Value = 0.1843806287874228
$ws.Range("O3").Value = 0.1919731617572378
$ws.Range("P3").Value = 0.1954407884571841
$ws.Range("S3").Value = 0.02514269694050685
$ws.Range("T3").Value = 0.03603549546644529

# Row 4 (Sending=FAPs, Target=Inflammatory-Mac)
$ws.Range("I4").Value = 0.1309698538606213
$ws.Range("J4").Value = 0.1843806287874228
$ws.Range("M4").Value = 2.312753666666667
$ws.Range("N4").Value = 6.938261000000001
$ws.Range("O4").Value = 0.4488476538468058
$ws.Range("P4").Value = 0.4569552252095968
$ws.Range("Q4").Value = 0.7416415111404445
$ws.Range("R4").Value = 6.674773600264
$ws.Range("S4").Value = 0.0587855116299989
$ws.Range("T4").Value = 0.08425369175184387

# Row 5 (Sending=FAPs, Target=MuSCs)
$ws.Range("I5").Value = 0.1309698538606213
$ws.Range("J5").Value = 0.1843806287874228
$ws.Range("M5").Value = 0.274264
$ws.Range("N5").Value = 0.548528
$ws.Range("O5").Value = 0.05322778413840601
$ws.Range("P5").Value = 0.03612616126343038
$ws.Range("Q5").Value = 0.08794951677866666
$ws.Range("R5").Value = 0.527697100672
$ws.Range("S5").Value = 0.006971235109931732
$ws.Range("T5").Value = 0.006660964329427132

# Row 6 (Sending=FAPs, Target=Resolving-Mac)
$ws.Range("I6").Value = 0.1309698538606213
$ws.Range("J6").Value = 0.1843806287874228
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.344694
$ws.Range("N6").Value = 4.034082
$ws.Range("O6").Value = 0.2609714799033402
$ws.Range("P6").Value = 0.2656854287874124
$ws.Range("Q6").Value = 0.4312093002186666
$ws.Range("R6").Value = 3.880883701968
$ws.Range("S6").Value = 0.03417939658473054
$ws.Range("T6").Value = 0.04898724641947914

# Row 7 (Sending=MuSCs, Target=ECs)
$ws.Range("G7").Value = 2.127787
$ws.Range("H7").Value = 4.255574
$ws.Range("I7").Value = 0.8690301461393787
$ws.Range("J7").Value = 0.8156193712125771
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.2317656666666667
$ws.Range("N7").Value = 0.6952970000000001
$ws.Range("O7").Value = 0.04497992035421015
$ws.Range("P7").Value = 0.04579239628237638
$ws.Range("Q7").Value = 0.4931479725796667
$ws.Range("R7").Value = 2.958887835478
$ws.Range("S7").Value = 0.03908890675875686
$ws.Range("T7").Value = 0.03734916546214898

# Row 8 (Sending=MuSCs, Target=FAPs)
$ws.Range("G8").Value = 2.127787
$ws.Range("H8").Value = 4.255574
$ws.Range("I8").Value = 0.8690301461393787
$ws.Range("J8").Value = 0.8156193712125771
$ws.Range("O8").Value = 0.1919731617572378
$ws.Range("P8").Value = 0.1954407884571841
$ws.Range("Q8").Value = 2.10474306679
$ws.Range("R8").Value = 12.62845840074
$ws.Range("S8").Value = 0.1668304648167309
$ws.Range("T8").Value = 0.1594052929907388

# Row 9 (Sending=MuSCs, Target=Inflammatory-Mac)
$ws.Range("G9").Value = 2.127787
$ws.Range("H9").Value = 4.255574
$ws.Range("I9").Value = 0.8690301461393787
$ws.Range("J9").Value = 0.8156193712125771
$ws.Range("M9").Value = 2.312753666666667
$ws.Range("N9").Value = 6.938261000000001
$ws.Range("O9").Value = 0.4488476538468058
$ws.Range("P9").Value = 0.4569552252095968
$ws.Range("Q9").Value = 4.921047186135667
$ws.Range("R9").Value = 29.526283116814
$ws.Range("S9").Value = 0.3900621422168069
$ws.Range("T9").Value = 0.3727015334577529

# Row 10 (Sending=MuSCs, Target=MuSCs)
$ws.Range("G10").Value = 2.127787
$ws.Range("H10").Value = 4.255574
$ws.Range("I10").Value = 0.8690301461393787
$ws.Range("J10").Value = 0.8156193712125771
$ws.Range("M10").Value = 0.274264
$ws.Range("N10").Value = 0.548528
$ws.Range("O10").Value = 0.05322778413840601
$ws.Range("P10").Value = 0.03612616126343038
$ws.Range("Q10").Value = 0.5835753737680001
$ws.Range("R10").Value = 2.334301495072
$ws.Range("S10").Value = 0.04625654902847428
$ws.Range("T10").Value = 0.02946519693400325

# Row 11 (Sending=MuSCs, Target=Resolving-Mac)
$ws.Range("G11").Value = 2.127787
$ws.Range("H11").Value = 4.255574
$ws.Range("I11").Value = 0.8690301461393787
$ws.Range("J11").Value = 0.8156193712125771
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 1.344694
$ws.Range("N11").Value = 4.034082
$ws.Range("O11").Value = 0.2609714799033402
$ws.Range("P11").Value = 0.2656854287874124
$ws.Range("Q11").Value = 2.861222412178
$ws.Range("R11").Value = 17.167334473068
$ws.Range("S11").Value = 0.2267920833186097
$ws.Range("T11").Value = 0.2166981823679332
